$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment is safe
$ws.Range('D2').Value = '67.365.31'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '3.518.93'
$ws.Range('E3').Value = '  -1.03%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('E6').Value = '  -1.90%  '
$ws.Range('D7').Value = '3.517.13'
$ws.Range('E7').Value = '  -1.02%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('E10').Value = '  -0.93%  '
$ws.Range('E11').Value = '  +2.16%  '
$ws.Range('E12').Value = '  -1.34%  '
$ws.Range('E13').Value = '  -1.22%  '
$ws.Range('D14').Value = '4.111.20'
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('E15').Value = '  -0.70%  '
$ws.Range('D16').Value = '3.517.48'
$ws.Range('E16').Value = '  -0.87%  '
$ws.Range('D17').Value = '67.327.81'
$ws.Range('E17').Value = '  -0.51%  '
$ws.Range('E18').Value = '  +0.28%  '
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('E20').Value = '  -2.53%  '
$ws.Range('E21').Value = '  -2.27%  '
$ws.Range('E22').Value = '  -1.10%  '
$ws.Range('E23').Value = '  -2.19%  '
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('E25').Value = '  +9.85%  '
$ws.Range('D26').Value = '3.656.33'
$ws.Range('E26').Value = '  -1.07%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('E27').Value = '  -2.13%  '
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('E29').Value = '  +0.45%  '
$ws.Range('E30').Value = '  -2.49%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('E31').Value = '  -4.83%  '
$ws.Range('B32').Value = 'Binance-PegBSC-USD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('E33').Value = '  +3.46%  '
$ws.Range('E34').Value = '  -0.19%  '
$ws.Range('E35').Value = '  -1.12%  '
$ws.Range('D36').Value = '3.506.79'
$ws.Range('E36').Value = '  -1.26%  '
$ws.Range('E37').Value = '  -3.09%  '
$ws.Range('E38').Value = '  -0.54%  '
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E40').Value = '  +0.68%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('E42').Value = '  +4.57%  '
$ws.Range('E43').Value = '  -0.36%  '
$ws.Range('E44').Value = '  -2.76%  '
$ws.Range('E45').Value = '  -0.84%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('E46').Value = '  -1.96%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('E47').Value = '  -2.22%  '
$ws.Range('E48').Value = '  +1.88%  '
$ws.Range('E49').Value = '  +4.89%  '
$ws.Range('E50').Value = '  -1.20%  '
$ws.Range('E51').Value = '  -0.88%  '

# Numeric-looking values that must remain stored as text (matches original inlineStr formatting)
$numericTextCells = @(
    @{Cell='D4'; Value='0.999'},
    @{Cell='D5'; Value='611.67'},
    @{Cell='D6'; Value='151.33'},
    @{Cell='D10'; Value='0.140'},
    @{Cell='D11'; Value='7.04'},
    @{Cell='D12'; Value='0.426'},
    @{Cell='D13'; Value='0.0000222'},
    @{Cell='D15'; Value='31.96'},
    @{Cell='D21'; Value='444.02'},
    @{Cell='D22'; Value='9.36'},
    @{Cell='D24'; Value='77.32'},
    @{Cell='D27'; Value='10.35'},
    @{Cell='D28'; Value='1.00'},
    @{Cell='D29'; Value='8.39'},
    @{Cell='D31'; Value='1.55'},
    @{Cell='D32'; Value='1.00'},
    @{Cell='D33'; Value='0.164'},
    @{Cell='D34'; Value='25.90'},
    @{Cell='D35'; Value='6.17'},
    @{Cell='D38'; Value='8.02'},
    @{Cell='D40'; Value='178.00'},
    @{Cell='D41'; Value='1.00'},
    @{Cell='D43'; Value='0.0881'},
    @{Cell='D44'; Value='5.46'},
    @{Cell='D45'; Value='0.882'},
    @{Cell='D46'; Value='28.44'},
    @{Cell='D47'; Value='44.98'},
    @{Cell='D48'; Value='2.64'},
    @{Cell='D49'; Value='1.27'},
    @{Cell='D50'; Value='7.60'},
    @{Cell='D51'; Value='1.01'}
)

foreach ($item in $numericTextCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.Style = "Normal"
}

Write-Output "Applied $($numericTextCells.Count + 74) cell updates"
